# "fixed export and fixing maps"
#
# The sheet used to show three census years (1989 / 2002 / 2014) with a
# sub-title row referencing the population census. The edit:
#   - drops the "(according to the population census data)" sub-title
#     (its row is kept, but now blank)
#   - removes the now-unneeded blank spacer row underneath it
#   - keeps only the 2014 column of data (drops the 1989 and 2002 columns),
#     so the surviving "Area" value (504.24...) and its "2014" header end
#     up in column B
#   - the table rows end up with a uniform 20.1pt row height

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "(according to the population census data)" sub-title text;
# the row itself stays in place as a blank spacer.
$ws.Range("A2").ClearContents()

# Remove the (now redundant) blank spacer row that used to sit between the
# sub-title and the "(sq. km)" row - everything below shifts up one row.
$ws.Rows.Item(3).Delete()

# Drop the 1989 and 2002 columns entirely; the surviving 2014 column
# (old column D) shifts left into column B.
$ws.Columns.Item(2).Delete()
$ws.Columns.Item(2).Delete()

# The five remaining rows all get a uniform 20.1pt height.
$ws.Rows("1:5").RowHeight = 20.1
